# 19/12/2025: Update the list
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the completed UnifyApps / Darrell Van Dyken entry (old row 2).
# This shifts all subsequent rows up by one.
$ws.Rows.Item(2).Delete()

# Insert a new row for the Rox "ENT AE San Fran" opening (Daniel Rodriquez),
# placed right after the other two Rox rows (now rows 2 and 3).
$ws.Rows.Item(4).Insert()
$ws.Cells.Item(4, 1).Value = 724
$ws.Cells.Item(4, 2).Value = "Rox"
$ws.Cells.Item(4, 3).Value = "ENT AE San Fran"
$ws.Cells.Item(4, 4).Value = "Daniel Rodriquez"
$ws.Cells.Item(4, 5).Value = "CV Sent"

# Swap the LanceDB candidate order: Rob Zhu now comes before Stephen Burd.
$ws.Cells.Item(11, 4).Value = "Rob Zhu"
$ws.Cells.Item(11, 5).Value = "CV Sent"
$ws.Cells.Item(12, 4).Value = "Stephen Burd"
$ws.Cells.Item(12, 5).Value = "1st Interview"

# Gino Lucia progressed from 3rd to 4th interview.
$ws.Cells.Item(13, 5).Value = "4th Interview"

# GRANT CUOMO (Dash0) progressed from CV Sent to 1st Interview.
$ws.Cells.Item(17, 5).Value = "1st Interview"

# Tim Duffy (Dash0) progressed from CV Sent to 1st Interview.
$ws.Cells.Item(21, 5).Value = "1st Interview"

# Append two new candidates at the bottom of the list.
$ws.Cells.Item(24, 1).Value = 859
$ws.Cells.Item(24, 2).Value = "Orca Ai"
$ws.Cells.Item(24, 3).Value = "Director of APAC"
$ws.Cells.Item(24, 4).Value = "Vivaan Seth"
$ws.Cells.Item(24, 5).Value = "CV Sent"

$ws.Cells.Item(25, 1).Value = 866
$ws.Cells.Item(25, 2).Value = "CyCognito"
$ws.Cells.Item(25, 3).Value = "BDR US"
$ws.Cells.Item(25, 4).Value = "Tyler Drago"
$ws.Cells.Item(25, 5).Value = "CV Sent"
